$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix rows 17/18: gene rows got re-ordered (DDB1 <-> RNF113A) ---
# Row 17 becomes RNF113A / ring finger protein 113A / confidence 2
$ws.Range("B17").Value = "RNF113A"
$ws.Range("C17").Value = "ring finger protein 113A"
# Force the numeric-looking confidence value to stay text (matches source data),
# then drop back to the default (unstyled) cell style like the rest of column D.
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2"
$ws.Range("D17").Style = "Normal"

# Row 18 becomes DDB1 / damage specific DNA binding protein 1 / confidence 1
$ws.Range("B18").Value = "DDB1"
$ws.Range("C18").Value = "damage specific DNA binding protein 1"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1"
$ws.Range("D18").Style = "Normal"

# --- Add new "time_taken" column F ---
# Header cell: copy formatting (bold/border/alignment) from E1's header style
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$ws.Range("F2").Value = "2021-10-05 13:42:44.420720"
$ws.Range("F3").Value = "2021-10-05 13:42:44.420732"
$ws.Range("F4").Value = "2021-10-05 13:42:44.420735"
$ws.Range("F5").Value = "2021-10-05 13:42:44.420738"
$ws.Range("F6").Value = "2021-10-05 13:42:44.420741"
$ws.Range("F7").Value = "2021-10-05 13:42:44.420744"
$ws.Range("F8").Value = "2021-10-05 13:42:44.420747"
$ws.Range("F9").Value = "2021-10-05 13:42:44.420749"
$ws.Range("F10").Value = "2021-10-05 13:42:44.420752"
$ws.Range("F11").Value = "2021-10-05 13:42:44.420755"
$ws.Range("F12").Value = "2021-10-05 13:42:44.420758"
$ws.Range("F13").Value = "2021-10-05 13:42:44.420761"
$ws.Range("F14").Value = "2021-10-05 13:42:44.420764"
$ws.Range("F15").Value = "2021-10-05 13:42:44.420766"
$ws.Range("F16").Value = "2021-10-05 13:42:44.420769"
$ws.Range("F17").Value = "2021-10-05 13:42:44.420772"
$ws.Range("F18").Value = "2021-10-05 13:42:44.420775"
